# Apply updated crypto price/volume data as produced by the scheduled scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.961.48"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.360.43"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.681"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.37"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.27%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.98%  "
$ws.Range("E13").Value = "  +10.10%  "
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "2.711.69"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "2.362.66"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "43.911.49"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +25.97%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0770"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.208"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.62%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.66%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "58.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.35%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.31%  "
